# Update market-price-derived columns (H-N) across several leve-profit sheets.
# Values sourced from an external market data refresh (see commit message).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 4490.6206
$ws.Range("I70").Value = 1710.2222
$ws.Range("J70").Value = 5741.8
$ws.Range("K70").Value = 5130.6666
$ws.Range("L70").Value = 17225.4
$ws.Range("M70").Value = -4860.6666
$ws.Range("N70").Value = -17765.4

$ws.Range("H73").Value = 4490.6206
$ws.Range("I73").Value = 1710.2222
$ws.Range("J73").Value = 5741.8
$ws.Range("K73").Value = 5130.6666
$ws.Range("L73").Value = 17225.4
$ws.Range("M73").Value = -4194.6666
$ws.Range("N73").Value = -19097.4

$ws.Range("H86").Value = 2999.5
$ws.Range("I86").Value = 996.5
$ws.Range("K86").Value = 996.5
$ws.Range("M86").Value = 126.5

$ws.Range("H87").Value = 70118
$ws.Range("J87").Value = 70118
$ws.Range("L87").Value = 70118
$ws.Range("N87").Value = -72614

$ws.Range("H89").Value = 2999.5
$ws.Range("I89").Value = 996.5
$ws.Range("K89").Value = 4982.5
$ws.Range("M89").Value = 633.5

$ws.Range("H90").Value = 70118
$ws.Range("J90").Value = 70118
$ws.Range("L90").Value = 210354
$ws.Range("N90").Value = -222834

$ws.Range("H138").Value = 442680.3
$ws.Range("J138").Value = 462575.3
$ws.Range("L138").Value = 1387725.9
$ws.Range("N138").Value = -1398005.9

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 10678.923
$ws.Range("I61").Value = 12980.667
$ws.Range("J61").Value = 5500
$ws.Range("K61").Value = 12980.667
$ws.Range("L61").Value = 5500
$ws.Range("M61").Value = -12768.667
$ws.Range("N61").Value = -5924

$ws.Range("H74").Value = 2244.0667
$ws.Range("I74").Value = 1612.88
$ws.Range("J74").Value = 5400
$ws.Range("K74").Value = 1612.88
$ws.Range("L74").Value = 5400
$ws.Range("M74").Value = -738.8800000000001
$ws.Range("N74").Value = -7148

$ws.Range("H77").Value = 2244.0667
$ws.Range("I77").Value = 1612.88
$ws.Range("J77").Value = 5400
$ws.Range("K77").Value = 8064.400000000001
$ws.Range("L77").Value = 27000
$ws.Range("M77").Value = -3696.400000000001
$ws.Range("N77").Value = -35736

$ws.Range("H97").Value = 9260.846
$ws.Range("I97").Value = 3868.375
$ws.Range("J97").Value = 17888.8
$ws.Range("K97").Value = 3868.375
$ws.Range("L97").Value = 17888.8
$ws.Range("M97").Value = -3372.375
$ws.Range("N97").Value = -18880.8

$ws.Range("H136").Value = 10678.923
$ws.Range("I136").Value = 12980.667
$ws.Range("J136").Value = 5500
$ws.Range("K136").Value = 38942.001
$ws.Range("L136").Value = 16500
$ws.Range("M136").Value = -36392.001
$ws.Range("N136").Value = -21600

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 404.9
$ws.Range("I22").Value = 462.375
$ws.Range("K22").Value = 462.375
$ws.Range("M22").Value = -289.375

$ws.Range("H107").Value = 4467.778
$ws.Range("I107").Value = 5171.143
$ws.Range("J107").Value = 2006
$ws.Range("K107").Value = 5171.143
$ws.Range("L107").Value = 2006
$ws.Range("M107").Value = -3251.143
$ws.Range("N107").Value = -5846

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3139.7
$ws.Range("I31").Value = 1783.9565
$ws.Range("J31").Value = 4294.593
$ws.Range("K31").Value = 1783.9565
$ws.Range("L31").Value = 4294.593
$ws.Range("M31").Value = -1488.9565
$ws.Range("N31").Value = -4884.593

$ws.Range("H34").Value = 3139.7
$ws.Range("I34").Value = 1783.9565
$ws.Range("J34").Value = 4294.593
$ws.Range("K34").Value = 1783.9565
$ws.Range("L34").Value = 4294.593
$ws.Range("M34").Value = -1581.9565
$ws.Range("N34").Value = -4698.593

$ws.Range("H107").Value = 921.5
$ws.Range("I107").Value = 783
$ws.Range("J107").Value = 1060
$ws.Range("K107").Value = 783
$ws.Range("L107").Value = 1060
$ws.Range("M107").Value = 1137
$ws.Range("N107").Value = -4900

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 194
$ws.Range("J2").Value = 158.6
$ws.Range("L2").Value = 951.5999999999999
$ws.Range("N2").Value = -1177.6

$ws.Range("H17").Value = 1735.7142
$ws.Range("I17").Value = 1383.3334
$ws.Range("J17").Value = 2000
$ws.Range("K17").Value = 4150.0002
$ws.Range("L17").Value = 6000
$ws.Range("M17").Value = -3981.0002
$ws.Range("N17").Value = -6338

$ws.Range("H34").Value = 1450.3636
$ws.Range("I34").Value = 1450.3636
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 4351.0908
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = -4267.0908
$ws.Range("N34").ClearContents()

$ws.Range("H39").Value = 14797
$ws.Range("I39").Value = 1320.5
$ws.Range("J39").Value = 41750
$ws.Range("K39").Value = 3961.5
$ws.Range("L39").Value = 125250
$ws.Range("M39").Value = -3667.5
$ws.Range("N39").Value = -125838

$ws.Range("H55").Value = 8031.8
$ws.Range("J55").Value = 9790
$ws.Range("L55").Value = 29370
$ws.Range("N55").Value = -29724

$ws.Range("H68").Value = 6122.88
$ws.Range("I68").Value = 2667.111
$ws.Range("J68").Value = 8066.75
$ws.Range("K68").Value = 8001.333
$ws.Range("L68").Value = 24200.25
$ws.Range("M68").Value = -7190.333
$ws.Range("N68").Value = -25822.25

$ws.Range("H71").Value = 6122.88
$ws.Range("I71").Value = 2667.111
$ws.Range("J71").Value = 8066.75
$ws.Range("K71").Value = 24003.999
$ws.Range("L71").Value = 72600.75
$ws.Range("M71").Value = -19947.999
$ws.Range("N71").Value = -80712.75

$ws.Range("H107").Value = 590.5833
$ws.Range("I107").Value = 252
$ws.Range("J107").Value = 832.4286
$ws.Range("K107").Value = 756
$ws.Range("L107").Value = 2497.2858
$ws.Range("M107").Value = 1164
$ws.Range("N107").Value = -6337.2858

$ws.Range("H132").Value = 2040.2667
$ws.Range("J132").Value = 2171.8333
$ws.Range("L132").Value = 19546.4997
$ws.Range("N132").Value = -24606.4997

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 737.6923
$ws.Range("I107").Value = 460.4
$ws.Range("J107").Value = 911
$ws.Range("K107").Value = 460.4
$ws.Range("L107").Value = 911
$ws.Range("M107").Value = 1459.6
$ws.Range("N107").Value = -4751
